$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(4)
$r = $p.Range
# Exclude the trailing paragraph mark from the range so we only delete the run contents.
$r.End = $r.End - 1
$r.Text = ""
